$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = '@'
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue 'D2' '67.265.06'
Set-TextValue 'E2' '  +0.89%  '
Set-TextValue 'D3' '3.515.95'
Set-TextValue 'E3' '  +0.42%  '
Set-TextValue 'E4' '  +0.01%  '
Set-TextValue 'D5' '596.60'
Set-TextValue 'E5' '  +0.92%  '
Set-TextValue 'D6' '173.55'
Set-TextValue 'E6' '  +1.69%  '
Set-TextValue 'E7' '  +0.03%  '
Set-TextValue 'D8' '0.594'
Set-TextValue 'E8' '  +2.53%  '
Set-TextValue 'E9' '  +5.82%  '
Set-TextValue 'E10' '  -0.60%  '
Set-TextValue 'D11' '0.437'
Set-TextValue 'E11' '  +0.04%  '
Set-TextValue 'D12' '4.127.30'
Set-TextValue 'E13' '  +0.03%  '
Set-TextValue 'D14' '29.30'
Set-TextValue 'E14' '  +3.59%  '
Set-TextValue 'D15' '67.201.38'
Set-TextValue 'E15' '  +0.77%  '
Set-TextValue 'D16' '0.0000180'
Set-TextValue 'E16' '  +1.17%  '
Set-TextValue 'D17' '3.509.34'
Set-TextValue 'E17' '  +0.37%  '
Set-TextValue 'D18' '6.36'
Set-TextValue 'E18' '  +0.29%  '
Set-TextValue 'D19' '14.19'
Set-TextValue 'E19' '  +1.28%  '
Set-TextValue 'D20' '395.86'
Set-TextValue 'E20' '  +1.96%  '
Set-TextValue 'D21' '8.04'
Set-TextValue 'E21' '  +0.46%  '
Set-TextValue 'D22' '73.17'
Set-TextValue 'E22' '  +0.13%  '
Set-TextValue 'E23' '  +0.11%  '
Set-TextValue 'D24' '0.539'
Set-TextValue 'E24' '  +1.03%  '
Set-TextValue 'D25' '0.0000122'
Set-TextValue 'E25' '  +0.53%  '
Set-TextValue 'D26' '10.31'
Set-TextValue 'E26' '  -1.39%  '
Set-TextValue 'E27' '  +1.10%  '
Set-TextValue 'E28' '  -0.23%  '
Set-TextValue 'D29' '6.34'
Set-TextValue 'E29' '  -0.14%  '
Set-TextValue 'D30' '1.47'
Set-TextValue 'E30' '  -0.55%  '
Set-TextValue 'E31' '  +0.34%  '
Set-TextValue 'D32' '23.87'
Set-TextValue 'E32' '  +1.30%  '
Set-TextValue 'D33' '7.40'
Set-TextValue 'E33' '  -0.57%  '
Set-TextValue 'D34' '1.69'
Set-TextValue 'E34' '  +4.20%  '
Set-TextValue 'D35' '162.92'
Set-TextValue 'E35' '  +0.17%  '
Set-TextValue 'E36' '  +1.02%  '
Set-TextValue 'E37' '  +0.64%  '
Set-TextValue 'D38' '7.04'
Set-TextValue 'E38' '  +5.88%  '
Set-TextValue 'E39' '  +0.33%  '
Set-TextValue 'E40' '  +0.77%  '
Set-TextValue 'D41' '26.68'
Set-TextValue 'E41' '  +0.91%  '
Set-TextValue 'D42' '27.26'
Set-TextValue 'E42' '  +2.52%  '
Set-TextValue 'D43' '2.843.83'
Set-TextValue 'E43' '  +0.77%  '
Set-TextValue 'D44' '2.59'
Set-TextValue 'E44' '  +2.91%  '
Set-TextValue 'D45' '43.01'
Set-TextValue 'E45' '  +0.04%  '
Set-TextValue 'D46' '0.0305'
Set-TextValue 'E46' '  -1.95%  '
Set-TextValue 'D47' '339.54'
Set-TextValue 'E47' '  -4.38%  '
Set-TextValue 'E48' '  +0.10%  '
Set-TextValue 'D49' '34.63'
Set-TextValue 'E49' '  +2.43%  '
Set-TextValue 'E50' '  +0.10%  '
Set-TextValue 'D51' '0.851'
Set-TextValue 'E51' '  -0.50%  '
